$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new status column G
$ws.Range("G1").Value = "Implementeret"

# Mark the already-implemented notes as "Done"
$ws.Range("G6").Value = "Done"
$ws.Range("G7").Value = "Done"
$ws.Range("G8").Value = "Done"
$ws.Range("G9").Value = "Done"
$ws.Range("G11").Value = "Done"

# Column widths (values chosen so the engine's pixel-quantized ColumnWidth
# setter lands on the exact stored widths used by the target workbook)
$ws.Columns.Item(2).ColumnWidth = 29.166666666666668
$ws.Columns.Item(6).ColumnWidth = 64.16666666666667
$ws.Columns.Item(7).ColumnWidth = 11.451822916666666

# Update the selection to reflect the new working cell
$ws.Range("F17").Select()
